# "extensão do ano máximo até 2060"
# Extend the "crescimento_mercado" (market growth) table from year 2050
# (row 42) down to year 2060 (row 52), repeating the 3% growth rate used
# for the preceding years, and leave the selection on the last-edited
# cell (B41) as recorded in the saved workbook view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow  = 43
$firstNewYear = 2051
$lastNewYear  = 2060
$growthRate   = 0.03

for ($year = $firstNewYear; $year -le $lastNewYear; $year++) {
    $row = $firstNewRow + ($year - $firstNewYear)
    $ws.Cells.Item($row, 1).Value = $year
    $ws.Cells.Item($row, 2).Value = $growthRate
}

# Restore the selection recorded in the workbook (cursor left on B41
# after the last edit, with the sheet scrolled down towards the new rows).
$ws.Range("B41").Select()
